$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new date labels for S1:U1, matching style of existing header cells
$ws.Range("S1").Value = "31/12/2023"
$ws.Range("T1").Value = "31/03/2024"
$ws.Range("U1").Value = "30/06/2024"
$ws.Range("R1").Copy()
$ws.Range("S1:U1").PasteSpecial(-4122)

# Data rows 2-80: new quarterly figures for columns S (31/12/2023), T (31/03/2024), U (30/06/2024)
$ws.Range("S2").Value = 3241590.016
$ws.Range("T2").Value = 3209953.536
$ws.Range("U2").Value = 3444075.52
$ws.Range("S3").Value = 2056580.096
$ws.Range("T3").Value = 2002643.968
$ws.Range("U3").Value = 2120319.872
$ws.Range("S4").Value = 221495.2
$ws.Range("T4").Value = 270202.976
$ws.Range("U4").Value = 274119.232
$ws.Range("S5").Value = 82994.52800000001
$ws.Range("T5").Value = 95284.696
$ws.Range("U5").Value = 67741.024
$ws.Range("S6").Value = 830832.3199999999
$ws.Range("T6").Value = 647678.464
$ws.Range("U6").Value = 734684.8639999999
$ws.Range("S7").Value = 782706.432
$ws.Range("T7").Value = 841014.912
$ws.Range("U7").Value = 884934.784
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = 0
$ws.Range("S9").Value = 118353.472
$ws.Range("T9").Value = 129255.864
$ws.Range("U9").Value = 139206.64
$ws.Range("S10").Value = 20198.152
$ws.Range("T10").Value = 19207.11
$ws.Range("U10").Value = 19633.346
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = 0
$ws.Range("U11").Value = 0
$ws.Range("S12").Value = 360073.088
$ws.Range("T12").Value = 356650.208
$ws.Range("U12").Value = 452394.144
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 0
$ws.Range("U13").Value = 0
$ws.Range("S14").Value = 0
$ws.Range("T14").Value = 0
$ws.Range("U14").Value = 0
$ws.Range("S15").Value = 0
$ws.Range("T15").Value = 0
$ws.Range("U15").Value = 0
$ws.Range("S16").Value = 0
$ws.Range("T16").Value = 0
$ws.Range("U16").Value = 0
$ws.Range("S17").Value = 0
$ws.Range("T17").Value = 0
$ws.Range("U17").Value = 0
$ws.Range("S18").Value = 0
$ws.Range("T18").Value = 0
$ws.Range("U18").Value = 0
$ws.Range("S19").Value = 222152.144
$ws.Range("T19").Value = 251474.112
$ws.Range("U19").Value = 267197.536
$ws.Range("S20").Value = 0
$ws.Range("T20").Value = 0
$ws.Range("U20").Value = 0
$ws.Range("S21").Value = 0
$ws.Range("T21").Value = 0
$ws.Range("U21").Value = 0
$ws.Range("S22").Value = 0
$ws.Range("T22").Value = 0
$ws.Range("U22").Value = 0
$ws.Range("S23").Value = 765745.9840000001
$ws.Range("T23").Value = 787747.776
$ws.Range("U23").Value = 798605.632
$ws.Range("S24").Value = 59190.864
$ws.Range("T24").Value = 62911.52
$ws.Range("U24").Value = 72755.728
$ws.Range("S25").Value = 0
$ws.Range("T25").Value = 0
$ws.Range("U25").Value = 0
$ws.Range("S26").Value = 3241590.016
$ws.Range("T26").Value = 3209953.536
$ws.Range("U26").Value = 3444075.52
$ws.Range("S27").Value = 695587.1360000001
$ws.Range("T27").Value = 677308.8639999999
$ws.Range("U27").Value = 692059.392
$ws.Range("S28").Value = 117858.768
$ws.Range("T28").Value = 96175.288
$ws.Range("U28").Value = 102690.208
$ws.Range("S29").Value = 103148.616
$ws.Range("T29").Value = 143560.784
$ws.Range("U29").Value = 210548.592
$ws.Range("S30").Value = 85080.82399999999
$ws.Range("T30").Value = 53753.632
$ws.Range("U30").Value = 82613.496
$ws.Range("S31").Value = 111463.432
$ws.Range("T31").Value = 126091.968
$ws.Range("U31").Value = 126079.048
$ws.Range("S32").Value = 0
$ws.Range("T32").Value = 0
$ws.Range("U32").Value = 0
$ws.Range("S33").Value = 87700.928
$ws.Range("T33").Value = 87700.928
$ws.Range("U33").Value = 6.343
$ws.Range("S34").Value = 190334.592
$ws.Range("T34").Value = 170026.256
$ws.Range("U34").Value = 170121.737
$ws.Range("S35").Value = 0
$ws.Range("T35").Value = 0
$ws.Range("U35").Value = 0
$ws.Range("S36").Value = 0
$ws.Range("T36").Value = 0
$ws.Range("U36").Value = 0
$ws.Range("S37").Value = 612280.8320000001
$ws.Range("T37").Value = 562700.032
$ws.Range("U37").Value = 578003.84
$ws.Range("S38").Value = 160000
$ws.Range("T38").Value = 100000
$ws.Range("U38").Value = 107558.992
$ws.Range("S39").Value = 0
$ws.Range("T39").Value = 0
$ws.Range("U39").Value = 0
$ws.Range("S40").Value = 437706.88
$ws.Range("T40").Value = 448765.472
$ws.Range("U40").Value = 455812.288
$ws.Range("S41").Value = 0
$ws.Range("T41").Value = 0
$ws.Range("U41").Value = 0
$ws.Range("S42").Value = 0
$ws.Range("T42").Value = 0
$ws.Range("U42").Value = 0
$ws.Range("S43").Value = 14573.945
$ws.Range("T43").Value = 13934.537
$ws.Range("U43").Value = 14632.574
$ws.Range("S44").Value = 0
$ws.Range("T44").Value = 0
$ws.Range("U44").Value = 0
$ws.Range("S45").Value = 0
$ws.Range("T45").Value = 0
$ws.Range("U45").Value = 0
$ws.Range("S46").Value = 0
$ws.Range("T46").Value = 0
$ws.Range("U46").Value = 0
$ws.Range("S47").Value = 1933722.112
$ws.Range("T47").Value = 1969944.576
$ws.Range("U47").Value = 2174012.16
$ws.Range("S48").Value = 1105381.248
$ws.Range("T48").Value = 1105381.248
$ws.Range("U48").Value = 1105381.248
$ws.Range("S49").Value = -68277.38400000001
$ws.Range("T49").Value = -67864.03200000001
$ws.Range("U49").Value = -74761.408
$ws.Range("S50").Value = 0
$ws.Range("T50").Value = 0
$ws.Range("U50").Value = 0
$ws.Range("S51").Value = 896618.24
$ws.Range("T51").Value = 896618.24
$ws.Range("U51").Value = 896619.7120000001
$ws.Range("S52").Value = 0
$ws.Range("T52").Value = 35809.184
$ws.Range("U52").Value = 246770.688
$ws.Range("S53").Value = 0
$ws.Range("T53").Value = 0
$ws.Range("U53").Value = 0
$ws.Range("S54").Value = 0
$ws.Range("T54").Value = 0
$ws.Range("U54").Value = 0
$ws.Range("S55").Value = 0
$ws.Range("T55").Value = 0
$ws.Range("U55").Value = 1.999
$ws.Range("S56").Value = 0
$ws.Range("T56").Value = 0
$ws.Range("U56").Value = 0
# Row 57 is a blank separator row in the source data; left unset
# Row 58 is a blank separator row in the source data; left unset
$ws.Range("S59").Value = 778127.232
$ws.Range("T59").Value = 444590.464
$ws.Range("U59").Value = 656322.048
$ws.Range("S60").Value = -233741.92
$ws.Range("T60").Value = -141600.928
$ws.Range("U60").Value = -204560.88
$ws.Range("S61").Value = 544385.28
$ws.Range("T61").Value = 302989.536
$ws.Range("U61").Value = 451761.152
$ws.Range("S62").Value = -276421.632
$ws.Range("T62").Value = -196293.776
$ws.Range("U62").Value = -245706.304
$ws.Range("S63").Value = -80130.32000000001
$ws.Range("T63").Value = -61157.68
$ws.Range("U63").Value = -63366.6
$ws.Range("S64").Value = 0
$ws.Range("T64").Value = 0
$ws.Range("U64").Value = 0
$ws.Range("S65").Value = 1667.35
$ws.Range("T65").Value = 1170.336
$ws.Range("U65").Value = 84353.152
$ws.Range("S66").Value = -542.86
$ws.Range("T66").Value = -5362.215
$ws.Range("U66").Value = -3141.774
$ws.Range("S67").Value = 0
$ws.Range("T67").Value = 0
$ws.Range("U67").Value = 0
$ws.Range("S68").Value = -17253.722
$ws.Range("T68").Value = -15499.737
$ws.Range("U68").Value = 3494.832
$ws.Range("S69").Value = 6415.464
$ws.Range("T69").Value = 7977.256
$ws.Range("U69").Value = 27911.692
$ws.Range("S70").Value = -23669.186
$ws.Range("T70").Value = -23476.992
$ws.Range("U70").Value = -24416.86
# Row 71 is a blank separator row in the source data; left unset
# Row 72 is a blank separator row in the source data; left unset
# Row 73 is a blank separator row in the source data; left unset
$ws.Range("S74").Value = 171704.064
$ws.Range("T74").Value = 25846.46
$ws.Range("U74").Value = 227394.448
$ws.Range("S75").Value = -3761.779
$ws.Range("T75").Value = -19359.248
$ws.Range("U75").Value = -32156.376
$ws.Range("S76").Value = -23792.068
$ws.Range("T76").Value = 29321.972
$ws.Range("U76").Value = 15723.418
# Row 77 is a blank separator row in the source data; left unset
# Row 78 is a blank separator row in the source data; left unset
$ws.Range("S79").Value = 0
$ws.Range("T79").Value = 0
$ws.Range("U79").Value = 0
$ws.Range("S80").Value = 144150.224
$ws.Range("T80").Value = 35809.184
$ws.Range("U80").Value = 210961.504
